$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $p.LineSpacingRule = 0   # wdLineSpaceSingle
}
